$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "stones" sub-table column (G) for the first few rows
$ws.Range("G1").Value2 = "stones"
$ws.Range("G2").Value2 = "agate (blue)"
$ws.Range("G3").Value2 = "nothing"
$ws.Range("G4").Value2 = "quartz (clear)"

# Fix "sandels" typo -> "sandals"
$ws.Range("F3").Value2 = "sandals"
$ws.Range("F10").Value2 = "sandals"
$ws.Range("F14").Value2 = "sandals"
$ws.Range("F18").Value2 = "sandals"

# New random item
$ws.Range("F16").Value2 = "rotten carrot"

# Re-used item for F8
$ws.Range("F8").Value2 = "rags"

# Update selection to reflect the cell used while editing the new column
$ws.Range("H13").Select() | Out-Null
